$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the 4 existing cell styles into scratch cells (col Z) so later
#     writes cannot disturb a style before it has been copied elsewhere. ---
$ws.Range("B1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Re-style the rows whose formatting differs from the new roster text. ---
$ws.Range("Z3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("Z4").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Remove the scratch cells; they must not remain in the saved sheet. ---
$ws.Range("Z1:Z4").Clear()

# --- Write the new roster list (column B, rows 2-32). ---
$ws.Range("B2").Value = "林健太郎"
$ws.Range("B3").Value = "山口洸翔"
$ws.Range("B4").Value = "志塚惇希"
$ws.Range("B5").Value = "小溝賢"
$ws.Range("B6").Value = "小野文哉"
$ws.Range("B7").Value = "渡部魁"
$ws.Range("B8").Value = "崎谷航平"
$ws.Range("B9").Value = ""
$ws.Range("B10").Value = "白岩詩佑介"
$ws.Range("B11").Value = "三神佳誠"
$ws.Range("B12").Value = "氏家琉貴"
$ws.Range("B13").Value = "羽賀尚生"
$ws.Range("B14").Value = "足立耕平"
$ws.Range("B15").Value = "遠藤隼人"
$ws.Range("B16").Value = "Ethan Virtudazo"
$ws.Range("B17").Value = "富澤天音"
$ws.Range("B18").Value = "Owen Williamson"
$ws.Range("B19").Value = "池田伊吹"
$ws.Range("B20").Value = "神山修造"
$ws.Range("B21").Value = "川田涼介"
$ws.Range("B22").Value = "豊島亮"
$ws.Range("B23").Value = "高野怜央"
$ws.Range("B24").Value = ""
$ws.Range("B25").Value = ""
$ws.Range("B26").Value = ""
$ws.Range("B27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("B30").Value = ""
$ws.Range("B31").Value = ""
$ws.Range("B32").Value = ""

# --- View: zoomed out from 171% to 75%, no pinned top-left row, selection -> G28. ---
$excel.ActiveWindow.Zoom = 75
$ws.Range("G28").Select()
